# Auto-update draw results: append the 2025-10-01 Pick 4 row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

# A15 and C15 look like a date / a plain number to Excel's auto-detection,
# so they need the "treat as text" quote-prefix to stay literal text
# (matching how all the other rows in this column are stored).
$ws.Range("A$row").Value = "'2025-10-01"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "'251001"
$ws.Range("D$row").Value = "8-6-4-0"
$ws.Range("E$row").Value = "2025-10-01T21:38:24.555+04:00"
